$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = Get-Date -Year 2016 -Month 9 -Day 12 -Hour 21 -Minute 15 -Second 35
$ws.Range("B17").Value = 14
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 32
$ws.Range("E17").Value = 77
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 7385
$ws.Range("H17").Value = 8745
$ws.Range("I17").Value = 878
$ws.Range("J17").Value = 183
$ws.Range("K17").Value = 91
$ws.Range("L17").Value = 7
$ws.Range("M17").Value = 2
$ws.Range("N17").Value = "Bag"
